# Add a "Justifications (if any)" column (column H) to the test-cases sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header for the new column.
$ws.Range("H1").Value = "Justifications (if any)"

# Placeholder "-" for each existing test-case row (rows 2-10).
$ws.Range("H2:H10").Value = "-"

# Match formatting of the neighbouring column G (header style + bordered body style).
$ws.Range("G1:G10").Copy() | Out-Null
$ws.Range("H1:H10").PasteSpecial(-4122) | Out-Null

# Widen the new column to comfortably fit its header text.
$ws.Columns.Item(8).ColumnWidth = 22.666666666666668

# Reflect the edit in the current selection, same as a user would leave it.
$ws.Range("H1:H10").Select() | Out-Null
